$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Value = 2.24
$ws.Range("T2").Value = 1.51
$ws.Range("U2").Value = 2.84
$ws.Range("F4").Value = 4.7
$ws.Range("G4").Value = 5.6
$ws.Range("H4").Value = 1.61
$ws.Range("I4").Value = 1.75
$ws.Range("G7").Value = 1.13
$ws.Range("J7").Value = 1.09
$ws.Range("N7").Value = 1.1
$ws.Range("P7").Value = 2.36
$ws.Range("Q7").Value = 1.24
$ws.Range("R7").Value = 2.04
$ws.Range("S7").Value = 1.51
$ws.Range("U7").Value = 1.44
$ws.Range("AN7").Value = 2.76
$ws.Range("F8").Value = 1.31
$ws.Range("H8").Value = 6
$ws.Range("J8").Value = 1.09
$ws.Range("K8").Value = 8.199999999999999
$ws.Range("O8").Value = 1.14
$ws.Range("P8").Value = 2.58
$ws.Range("Q8").Value = 1.43
$ws.Range("R8").Value = 1.74
$ws.Range("S8").Value = 1.91
$ws.Range("T8").Value = 1.79
$ws.Range("U8").Value = 2
$ws.Range("V8").Value = 1.08
$ws.Range("W8").Value = 3.25
$ws.Range("AN8").Value = 4.2
$ws.Range("G9").Value = 3.65
$ws.Range("I9").Value = 2.46
$ws.Range("J9").Value = 3.45
$ws.Range("K9").Value = 3.95
$ws.Range("N9").Value = 3.75
$ws.Range("O9").Value = 1.29
$ws.Range("P9").Value = 1.96
$ws.Range("Q9").Value = 1.84
$ws.Range("R9").Value = 1.38
$ws.Range("S9").Value = 3.1
$ws.Range("T9").Value = 1.69
$ws.Range("U9").Value = 2.18
$ws.Range("AM9").Value = 100
$ws.Range("P10").Value = 1.83
$ws.Range("Q10").Value = 1.47
$ws.Range("R10").Value = 1.4
$ws.Range("S10").Value = 2.66
$ws.Range("V10").Value = 1.46
$ws.Range("X10").Value = 990
$ws.Range("G11").Value = 2.18
$ws.Range("H11").Value = 3.35
$ws.Range("I11").Value = 3.9
$ws.Range("J11").Value = 3.6
$ws.Range("N11").Value = 5.5
$ws.Range("S11").Value = 2.3
$ws.Range("V11").Value = 1.35
$ws.Range("AB11").Value = 17.5
$ws.Range("AN11").Value = 9.6
$ws.Range("G12").Value = 2.68
$ws.Range("J12").Value = 3.15
$ws.Range("P12").Value = 1.73
$ws.Range("Q12").Value = 1.96
$ws.Range("S12").Value = 3.55
$ws.Range("F13").Value = 6.2
$ws.Range("L13").Value = 1.22
$ws.Range("N13").Value = 2.46
$ws.Range("P13").Value = 2.46
$ws.Range("R13").Value = 1.55
$ws.Range("F14").Value = 2.12
$ws.Range("H14").Value = 3.15
$ws.Range("J14").Value = 3.55
$ws.Range("K14").Value = 4.2
$ws.Range("Q14").Value = 1.86
$ws.Range("S14").Value = 3.2
$ws.Range("T14").Value = 1.72
$ws.Range("W14").Value = 1.73
$ws.Range("X14").Value = 19.5
$ws.Range("Y14").Value = 16
$ws.Range("Z14").Value = 980
$ws.Range("AD14").Value = 990
$ws.Range("AF14").Value = 980
$ws.Range("AG14").Value = 990
$ws.Range("AM14").Value = 110
$ws.Range("AN14").Value = 980
$ws.Range("F15").Value = 1.22
$ws.Range("I15").Value = 1000
$ws.Range("J15").Value = 6.4
$ws.Range("K15").Value = 950
$ws.Range("N15").Value = 2.68
$ws.Range("P15").Value = 2.68
$ws.Range("Q15").Value = 1.32
$ws.Range("S15").Value = 1.89
$ws.Range("T15").Value = 1.04
$ws.Range("U15").Value = 1.04
$ws.Range("AB15").Value = 990
$ws.Range("AF15").Value = 1000
$ws.Range("AG15").Value = 990
$ws.Range("AJ15").Value = 1000
$ws.Range("AK15").Value = 1000
$ws.Range("AL15").Value = 1000
$ws.Range("AN15").Value = 1000
$ws.Range("F16").Value = 2.02
$ws.Range("G16").Value = 2.16
$ws.Range("I16").Value = 4.1
$ws.Range("J16").Value = 3.75
$ws.Range("K16").Value = 4.3
$ws.Range("N16").Value = 4.1
$ws.Range("Q16").Value = 1.74
$ws.Range("S16").Value = 2.9
$ws.Range("T16").Value = 1.65
$ws.Range("V16").Value = 1.33
$ws.Range("W16").Value = 1.87
$ws.Range("X16").Value = 22
$ws.Range("Y16").Value = 21
$ws.Range("Z16").Value = 36
$ws.Range("AB16").Value = 13
$ws.Range("AH16").Value = 20
$ws.Range("AJ16").Value = 30
$ws.Range("AL16").Value = 40
$ws.Range("AM16").Value = 95
$ws.Range("AN16").Value = 15.5
$ws.Range("F17").Value = 5
$ws.Range("H18").Value = 2.8
$ws.Range("J18").Value = 3.2
$ws.Range("N18").Value = 3.25
$ws.Range("O18").Value = 1.35
$ws.Range("P18").Value = 1.78
$ws.Range("Q18").Value = 2.04
$ws.Range("R18").Value = 1.29
$ws.Range("S18").Value = 3.7
$ws.Range("U18").Value = 2.06
$ws.Range("X18").Value = 15.5
$ws.Range("Z18").Value = 24
$ws.Range("AA18").Value = 55
$ws.Range("AE18").Value = 980
$ws.Range("AF18").Value = 21
$ws.Range("AH18").Value = 21
$ws.Range("AJ18").Value = 50
$ws.Range("AK18").Value = 38
$ws.Range("AL18").Value = 55
$ws.Range("AM18").Value = 130
$ws.Range("AN18").Value = 34
$ws.Range("AO18").Value = 980
$ws.Range("F19").Value = 3.9
$ws.Range("I19").Value = 1.96
$ws.Range("J19").Value = 3.9
$ws.Range("K19").Value = 5.1
$ws.Range("Q19").Value = 1.65
$ws.Range("R19").Value = 1.44
$ws.Range("S19").Value = 2.38
$ws.Range("T19").Value = 1.48
$ws.Range("V19").Value = 2.04
$ws.Range("X19").Value = 29
$ws.Range("F20").Value = 5.8
$ws.Range("I20").Value = 1.59
$ws.Range("K20").Value = 5.6
$ws.Range("Q20").Value = 1.47
$ws.Range("T20").Value = 1.63
$ws.Range("U20").Value = 2.28
$ws.Range("X20").Value = 990
$ws.Range("Y20").Value = 990
$ws.Range("AI20").Value = 980
$ws.Range("F21").Value = 5.2
$ws.Range("O21").Value = 1.3
$ws.Range("S21").Value = 3.25
$ws.Range("T21").Value = 1.04
$ws.Range("U21").Value = 1.04
$ws.Range("Y21").Value = 8.800000000000001
$ws.Range("AH21").Value = 25
$ws.Range("AJ21").Value = 180
$ws.Range("AM21").Value = 150
$ws.Range("AN21").Value = 120
$ws.Range("I22").Value = 3.8
$ws.Range("P22").Value = 1.59
$ws.Range("F23").Value = 2.14
$ws.Range("K23").Value = 5.4
$ws.Range("R23").Value = 1.47
$ws.Range("S23").Value = 2.38
$ws.Range("I24").Value = 17.5
$ws.Range("J24").Value = 8.199999999999999
$ws.Range("P24").Value = 3.65
$ws.Range("Q24").Value = 1.32
$ws.Range("R24").Value = 2.04
$ws.Range("S24").Value = 1.78
$ws.Range("T24").Value = 1.74
$ws.Range("U24").Value = 2.04
$ws.Range("W24").Value = 5.2
$ws.Range("X24").Value = 55
$ws.Range("Z24").Value = 190
$ws.Range("AA24").Value = 820
$ws.Range("AB24").Value = 16.5
$ws.Range("AC24").Value = 21
$ws.Range("AF24").Value = 11.5
$ws.Range("AG24").Value = 12.5
$ws.Range("AH24").Value = 32
$ws.Range("AJ24").Value = 11.5
$ws.Range("AK24").Value = 12.5
$ws.Range("AL24").Value = 29
$ws.Range("AN24").Value = 3.05
$ws.Range("AO24").Value = 180
